# Adressage.docx edit:
# - Add a new note paragraph (surrounded by blank paragraphs) right after
#   the "... switch f1" paragraph in the Siege_sociale section.
# - Add one extra blank paragraph right after the paragraph that carries
#   the _GoBack bookmark.

$d = $word.ActiveDocument

# Locate the "switch f1" paragraph (end of the Siege_sociale f1 0/0 line)
# and append three new paragraphs right after it: blank, note text, blank.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*switch f1*") {
        $target = $p
    }
}

$r = $target.Range
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()
$r.InsertParagraphAfter()

# The middle of the three freshly inserted paragraphs gets the note text.
$target.Next().Next().Range.InsertAfter( `
    "Le ping de routeur a routeur e passant par le siege sociale ne parche pas")

# Locate the paragraph holding the _GoBack bookmark and add one blank
# paragraph right after it.
$bm = $d.Bookmarks.Item("_GoBack")
$bmPara = $bm.Range.Paragraphs.Item(1)
$bmPara.Range.InsertParagraphAfter()
